$d = $word.ActiveDocument

# Locate the paragraph that immediately follows "Deleted UI element from View
# Controller" (an otherwise-empty paragraph) by matching its w14:paraId, since
# that's stable regardless of any prior edits in this session.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.WordOpenXML -match 'w14:paraId="62281A33"') {
        $anchor = $para
        break
    }
}

# Insert three new paragraphs right after it: blank, the new note, blank.
$anchor.Range.InsertParagraphAfter()
$blank1 = $d.Paragraphs($anchor.Index + 1)

$blank1.Range.InsertParagraphAfter()
$noteParagraph = $d.Paragraphs($blank1.Index + 1)
$noteParagraph.Range.Text = "Deleted most of previously commented out lines."

$noteParagraph.Range.InsertParagraphAfter()

# Remove the now-redundant blank paragraph that used to sit right after the
# "_GoBack" bookmark paragraph.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.WordOpenXML -match 'w14:paraId="1E176B3C"') {
        $para.Range.Delete()
        break
    }
}
